$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.944.43'
$ws.Range("E2").Value = '  +0.65%  '
$ws.Range("D3").Value = '2.294.11'
$ws.Range("E3").Value = '  +0.44%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '299.66'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.48%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '97.05'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.48%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.504'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.55%  '
$ws.Range("E8").Value = '  -0.04%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.502'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.36%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '33.66'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.28%  '
$ws.Range("E11").Value = '  +0.42%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '49.16'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -3.18%  '
$ws.Range("E13").Value = '  +2.93%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '17.09'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +11.78%  '
$ws.Range("E15").Value = '  +1.76%  '
$ws.Range("D16").Value = '2.646.39'
$ws.Range("E16").Value = '  +0.38%  '
$ws.Range("D17").Value = '2.302.02'
$ws.Range("E17").Value = '  +0.92%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.805'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.87%  '
$ws.Range("D19").Value = '42.888.69'
$ws.Range("E19").Value = '  +0.84%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.66'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.53%  '
$ws.Range("D21").Value = '0.0₃0901'
$ws.Range("E21").Value = '  +0.61%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.04'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.71%  '
$ws.Range("E23").Value = '  +1.20%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '236.36'
$ws.Range("D24").Style = "Normal"
$ws.Range("E25").Value = '  +5.14%  '
$ws.Range("E26").Value = '  -0.04%  '
$ws.Range("E27").Value = '  -1.53%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '24.27'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.08%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.18'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -5.26%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '166.39'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.02%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '33.71'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.05%  '
$ws.Range("E33").Value = '  +0.15%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.95'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.40%  '
$ws.Range("E35").Value = '  +5.87%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.42'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.56%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '16.75'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +3.31%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0696'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.25%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.83'
$ws.Range("D39").Style = "Normal"
$ws.Range("E40").Value = '  +0.98%  '
$ws.Range("E41").Value = '  -0.27%  '
$ws.Range("E42").Value = '  -0.13%  '
$ws.Range("E43").Value = '  -1.26%  '
$ws.Range("D44").Value = '1.985.04'
$ws.Range("E44").Value = '  +1.28%  '
$ws.Range("E45").Value = '  +0.64%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '9.88'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.14%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '17.54'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.26%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.83'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.36%  '
$ws.Range("D49").Value = '2.528.00'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '53.05'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.82%  '
$ws.Range("E51").Value = '  -2.26%  '
